# Add 2022-Q3 data
# ------------------------------------------------------------------
# 1) Insert a brand-new worksheet named "2022-Q3" right after "总计"
#    (i.e. before the existing "2021-Q4" sheet). All later sheets
#    keep their content, just shifted one tab to the right.
# 2) Populate the new sheet with the 2022-Q3 fund holding data, using
#    the same layout/style as the other quarterly sheets.
# 3) Update the "总计" (summary) sheet: insert a new row right under
#    the header for "2022-Q3" and push the existing rows down.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)          # "总计"
$q4_2021 = $wb.Worksheets.Item(2)           # "2021-Q4" (about to become #3)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)   # currently-active last tab

# ------------------------------------------------------------------
# 1) Create the new sheet before "2021-Q4" and name it
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($q4_2021)
$newSheet.Name = "2022-Q3"

# "2021-Q4" is now one tab further to the right; re-resolve it so we
# can borrow its formatting for the new sheet.
$q4_2021 = $wb.Worksheets.Item(3)

# Copy the header + a data-row's formatting onto the new sheet so it
# matches the look of the other quarterly sheets (bold/bordered header
# in B1:H1, bold/bordered index column in A2:A5).
$q4_2021.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)

$q4_2021.Range("A2:H2").Copy()
$newSheet.Range("A2:H5").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Fill in the 2022-Q3 data
# ------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B (fund code, e.g. "013413" keeps its leading zero) and D:G
# hold text values in every quarterly sheet (not numbers), so force
# Text format before writing them.
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:G5").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "501030"
$newSheet.Range("C2").Value = "汇添富中证环境治理指数（LOF）A"
$newSheet.Range("D2").Value = "3.12"
$newSheet.Range("E2").Value = "92.74"
$newSheet.Range("F2").Value = "2.07"
$newSheet.Range("G2").Value = "0.0646"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "164908"
$newSheet.Range("C3").Value = "交银施罗德中证环境治理指数（LOF）"
$newSheet.Range("D3").Value = "1.57"
$newSheet.Range("E3").Value = "93.62"
$newSheet.Range("F3").Value = "2.05"
$newSheet.Range("G3").Value = "0.0322"
$newSheet.Range("H3").Value = 9

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "501031"
$newSheet.Range("C4").Value = "汇添富中证环境治理指数（LOF）C"
$newSheet.Range("D4").Value = "1.30"
$newSheet.Range("E4").Value = "92.74"
$newSheet.Range("F4").Value = "2.07"
$newSheet.Range("G4").Value = "0.0269"
$newSheet.Range("H4").Value = 10

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "013413"
$newSheet.Range("C5").Value = "交银施罗德中证环境治理指数（LOF）C"
$newSheet.Range("D5").Value = "0.09"
$newSheet.Range("E5").Value = "93.62"
$newSheet.Range("F5").Value = "2.05"
$newSheet.Range("G5").Value = "0.0018"
$newSheet.Range("H5").Value = 9

# ------------------------------------------------------------------
# 3) Update the "总计" summary sheet: insert a row for 2022-Q3 right
#    after the header, pushing the other quarters down by one row.
# ------------------------------------------------------------------
$summary.Rows.Item(2).Insert()

# Match the index-column (A) formatting used by the rest of the rows.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.13

# Keep the running 0-based index column (A) sequential for every row.
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ------------------------------------------------------------------
# Restore the originally-active tab (it shifted one position right).
# ------------------------------------------------------------------
$lastSheet.Activate()
